$wb = $excel.ActiveWorkbook

# ---- Step 1: insert new sheet '2022-Q1' right before the '总计' sheet ----
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = '2022-Q1'
$q1.StandardWidth = 8
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# helper-less inline formatting: bold + thin border + center/top alignment
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}
function Set-TextFormat($cell) {
    $cell.NumberFormat = '@'
}

# ---- header row ----
$c = $q1.Cells.Item(1, 2); $c.Value = '基金代码'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 3); $c.Value = '基金名称'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 4); $c.Value = '基金规模'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 5); $c.Value = '股票总仓位'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 6); $c.Value = '仓位占比'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 7); $c.Value = '持有市值(亿元)'; Set-HeaderStyle $c
$c = $q1.Cells.Item(1, 8); $c.Value = '仓位排名'; Set-HeaderStyle $c

# ---- data rows ----
$r = 2
$c = $q1.Cells.Item($r, 1); $c.Value = 0; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '003293'
$q1.Cells.Item($r, 3).Value = '易方达科瑞灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '34.67'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '78.17'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.06'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '1.4076'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 3
$r = 3
$c = $q1.Cells.Item($r, 1); $c.Value = 1; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '519019'
$q1.Cells.Item($r, 3).Value = '大成景阳领先混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '10.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '92.80'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '7.56'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.7809'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 4
$c = $q1.Cells.Item($r, 1); $c.Value = 2; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '110012'
$q1.Cells.Item($r, 3).Value = '易方达科汇灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '15.73'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '75.64'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.04'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.6355'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 2
$r = 5
$c = $q1.Cells.Item($r, 1); $c.Value = 3; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '010673'
$q1.Cells.Item($r, 3).Value = '兴全中证800六个月持有期指数增强A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '19.06'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '97.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '2.92'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.5566'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 6
$r = 6
$c = $q1.Cells.Item($r, 1); $c.Value = 4; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '008381'
$q1.Cells.Item($r, 3).Value = '前海开源新兴产业混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '6.18'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '93.63'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '7.70'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.4759'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 7
$c = $q1.Cells.Item($r, 1); $c.Value = 5; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '164403'
$q1.Cells.Item($r, 3).Value = '前海开源沪港深农业主题精选灵活配置混合（LOF）'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '6.21'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '92.31'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '7.43'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.4614'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 8
$r = 8
$c = $q1.Cells.Item($r, 1); $c.Value = 6; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '010389'
$q1.Cells.Item($r, 3).Value = '易方达科益混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '7.10'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '92.94'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '6.34'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.4501'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 9
$c = $q1.Cells.Item($r, 1); $c.Value = 7; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '003857'
$q1.Cells.Item($r, 3).Value = '前海开源周期优选灵活配置混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '5.52'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.41'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '7.80'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.4306'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 10
$c = $q1.Cells.Item($r, 1); $c.Value = 8; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '090016'
$q1.Cells.Item($r, 3).Value = '大成消费主题混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '4.23'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '93.78'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '9.39'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.3972'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 2
$r = 11
$c = $q1.Cells.Item($r, 1); $c.Value = 9; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011649'
$q1.Cells.Item($r, 3).Value = '易方达逆向投资混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '7.49'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '85.02'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.55'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.3408'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 3
$r = 12
$c = $q1.Cells.Item($r, 1); $c.Value = 10; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '001320'
$q1.Cells.Item($r, 3).Value = '工银瑞信丰盈回报灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '6.50'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '93.76'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.14'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.2691'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 8
$r = 13
$c = $q1.Cells.Item($r, 1); $c.Value = 11; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '003858'
$q1.Cells.Item($r, 3).Value = '前海开源周期优选灵活配置混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '2.32'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.41'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '7.80'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.1810'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 14
$c = $q1.Cells.Item($r, 1); $c.Value = 12; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011911'
$q1.Cells.Item($r, 3).Value = '华夏消费优选混合型证券投资基金A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '7.18'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '82.18'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '2.35'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.1687'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 10
$r = 15
$c = $q1.Cells.Item($r, 1); $c.Value = 13; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '519756'
$q1.Cells.Item($r, 3).Value = '交银施罗德国企改革灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '5.69'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '87.70'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '2.77'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.1576'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 8
$r = 16
$c = $q1.Cells.Item($r, 1); $c.Value = 14; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '001907'
$q1.Cells.Item($r, 3).Value = '国投瑞银境煊灵活配置混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '2.61'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.67'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.1219'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 7
$r = 17
$c = $q1.Cells.Item($r, 1); $c.Value = 15; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011650'
$q1.Cells.Item($r, 3).Value = '易方达逆向投资混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.96'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '85.02'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.55'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0892'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 3
$r = 18
$c = $q1.Cells.Item($r, 1); $c.Value = 16; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '001908'
$q1.Cells.Item($r, 3).Value = '国投瑞银境煊灵活配置混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.75'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.67'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0817'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 7
$r = 19
$c = $q1.Cells.Item($r, 1); $c.Value = 17; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011431'
$q1.Cells.Item($r, 3).Value = '泰达宏利消费服务混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.61'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '81.15'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '3.19'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0514'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 20
$c = $q1.Cells.Item($r, 1); $c.Value = 18; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '000916'
$q1.Cells.Item($r, 3).Value = '前海开源股息率100强等权重股票'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '3.07'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.96'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.58'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0485'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 6
$r = 21
$c = $q1.Cells.Item($r, 1); $c.Value = 19; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '010674'
$q1.Cells.Item($r, 3).Value = '兴全中证800六个月持有期指数增强C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.49'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '97.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '2.92'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0435'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 6
$r = 22
$c = $q1.Cells.Item($r, 1); $c.Value = 20; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '002319'
$q1.Cells.Item($r, 3).Value = '大成一带一路灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.50'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '89.30'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '8.16'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0408'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 23
$c = $q1.Cells.Item($r, 1); $c.Value = 21; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '009726'
$q1.Cells.Item($r, 3).Value = '招商中证500等权重指数增强A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.87'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '91.11'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0269'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 10
$r = 24
$c = $q1.Cells.Item($r, 1); $c.Value = 22; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '008846'
$q1.Cells.Item($r, 3).Value = '大成民稳增长混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '2.49'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.89'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.04'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0259'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 9
$r = 25
$c = $q1.Cells.Item($r, 1); $c.Value = 23; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '010390'
$q1.Cells.Item($r, 3).Value = '易方达科益混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.29'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '92.94'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '6.34'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0184'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 26
$c = $q1.Cells.Item($r, 1); $c.Value = 24; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '003182'
$q1.Cells.Item($r, 3).Value = '华富弘鑫灵活配置混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '5.52'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.56'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '0.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0182'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 6
$r = 27
$c = $q1.Cells.Item($r, 1); $c.Value = 25; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '008037'
$q1.Cells.Item($r, 3).Value = '兴银先锋成长混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.41'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '79.32'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.36'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0179'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 1
$r = 28
$c = $q1.Cells.Item($r, 1); $c.Value = 26; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '015309'
$q1.Cells.Item($r, 3).Value = '国投瑞银境煊灵活配置混合E'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '90.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.67'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0154'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 7
$r = 29
$c = $q1.Cells.Item($r, 1); $c.Value = 27; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011912'
$q1.Cells.Item($r, 3).Value = '华夏消费优选混合型证券投资基金C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '82.18'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '2.35'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0103'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 10
$r = 30
$c = $q1.Cells.Item($r, 1); $c.Value = 28; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '009727'
$q1.Cells.Item($r, 3).Value = '招商中证500等权重指数增强C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.69'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '91.11'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.44'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0099'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 10
$r = 31
$c = $q1.Cells.Item($r, 1); $c.Value = 29; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '008038'
$q1.Cells.Item($r, 3).Value = '兴银先锋成长混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.17'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '79.32'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '4.36'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0074'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 1
$r = 32
$c = $q1.Cells.Item($r, 1); $c.Value = 30; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '008847'
$q1.Cells.Item($r, 3).Value = '大成民稳增长混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.56'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.89'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.04'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0058'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 9
$r = 33
$c = $q1.Cells.Item($r, 1); $c.Value = 31; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '009796'
$q1.Cells.Item($r, 3).Value = '大成汇享一年持有期混合A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.38'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.99'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.37'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0052'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 5
$r = 34
$c = $q1.Cells.Item($r, 1); $c.Value = 32; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '003183'
$q1.Cells.Item($r, 3).Value = '华富弘鑫灵活配置混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '1.54'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.56'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '0.33'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0051'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 6
$r = 35
$c = $q1.Cells.Item($r, 1); $c.Value = 33; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '011432'
$q1.Cells.Item($r, 3).Value = '泰达宏利消费服务混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.15'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '81.15'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '3.19'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0048'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 4
$r = 36
$c = $q1.Cells.Item($r, 1); $c.Value = 34; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '001474'
$q1.Cells.Item($r, 3).Value = '兴银丰盈灵活配置混合'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.08'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '83.45'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '3.59'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0029'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 2
$r = 37
$c = $q1.Cells.Item($r, 1); $c.Value = 35; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '004790'
$q1.Cells.Item($r, 3).Value = '富荣中证500指数增强A'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.06'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '84.00'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.49'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0009'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 5
$r = 38
$c = $q1.Cells.Item($r, 1); $c.Value = 36; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '004791'
$q1.Cells.Item($r, 3).Value = '富荣中证500指数增强C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.04'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '84.00'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.49'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0006'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 5
$r = 39
$c = $q1.Cells.Item($r, 1); $c.Value = 37; Set-HeaderStyle $c
$q1.Cells.Item($r, 2).Value = '009797'
$q1.Cells.Item($r, 3).Value = '大成汇享一年持有期混合C'
$c = $q1.Cells.Item($r, 4); Set-TextFormat $c; $c.Value = '0.04'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 5); Set-TextFormat $c; $c.Value = '22.99'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 6); Set-TextFormat $c; $c.Value = '1.37'; $c.ClearFormats()
$c = $q1.Cells.Item($r, 7); Set-TextFormat $c; $c.Value = '0.0005'; $c.ClearFormats()
$q1.Cells.Item($r, 8).Value = 5

# ---- Step 2: insert new summary row into '总计' sheet ----
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()
$c = $total.Cells.Item(2, 1); $c.Value = 0; Set-HeaderStyle $c
$total.Cells.Item(2, 2).Value = '2022-Q1'
$total.Cells.Item(2, 3).Value = 38
$total.Cells.Item(2, 4).Value = 7.37
